$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows 12, 13 and 14 need to be cyclically rotated:
#   new row12 = old row13
#   new row13 = old row14
#   new row14 = old row12
# Only columns A, B, E, F, G, H, P, Q, R, S change; all other columns
# (C, D, I, J, K, N, T, U, V, W, Y, AA, AD, AE, AF, AG, AT, AW, AX, AY, ...)
# are identical across the three rows and therefore remain untouched.

$cols = @("A", "B", "E", "F", "G", "H", "P", "Q", "R", "S")

# Capture current (before-edit) values for the three rows.
$row12 = @{}
$row13 = @{}
$row14 = @{}
foreach ($col in $cols) {
    $row12[$col] = $ws.Range("$col`12").Value2
    $row13[$col] = $ws.Range("$col`13").Value2
    $row14[$col] = $ws.Range("$col`14").Value2
}

# Apply the rotation.
foreach ($col in $cols) {
    $ws.Range("$col`12").Value2 = $row13[$col]
    $ws.Range("$col`13").Value2 = $row14[$col]
    $ws.Range("$col`14").Value2 = $row12[$col]
}
